$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cell D1
$ws.Range("D1").Value = "Measured Parameters"

# Row 2 additions: C2, D2 ; row height to 75
$ws.Range("C2").Value = "solid-state fermentation of the fungus Penicillium restrictum"
$ws.Range("D2").Value = "COD/Oil and Grease/Biogas/Methane Production/Free fatty acids (titrimetric)/ VS/ initial VSS"
$ws.Rows.Item(2).RowHeight = 75

# New row 3
$ws.Range("A3").Value = "Anaerobic digestion of lipid-rich waste — Effects of lipid concentration"
$ws.Range("B3").Value = 2006
$ws.Range("C3").Value = "Lipase 80,000 from Rhizopus oryzae"
$ws.Range("D3").Value = "Methane Production/GC-VFA/ VS/initial VSS"
$ws.Rows.Item(3).RowHeight = 45

# Column widths (target stored widths are 22.85546875 / 22.5703125 characters;
# the ColumnWidth setter internally rounds to 1/6-character steps, so we feed
# it the input value that lands closest to each exact target after rounding)
$ws.Columns.Item(3).ColumnWidth = 22
$ws.Columns.Item(4).ColumnWidth = 21.65

# Update selection to match target
$ws.Range("C4").Select()
